$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (29 Jun 2021 - 25 Jul 2021), appended after existing row 301.
$newData = @(
    @(44376, 0, 0, 0),
    @(44377, 0, 0, 0),
    @(44378, 0, 0, 0),
    @(44379, 0, 0, 0),
    @(44380, 1, 1, 26.76659528907923),
    @(44381, 0, 1, 26.76659528907923),
    @(44382, 0, 1, 26.76659528907923),
    @(44383, 0, 1, 26.76659528907923),
    @(44384, 0, 1, 26.76659528907923),
    @(44385, 0, 1, 26.76659528907923),
    @(44386, 0, 1, 26.76659528907923),
    @(44387, 0, 0, 0),
    @(44388, 0, 0, 0),
    @(44389, 0, 0, 0),
    @(44390, 0, 0, 0),
    @(44391, 0, 0, 0),
    @(44392, 0, 0, 0),
    @(44393, 0, 0, 0),
    @(44394, 0, 0, 0),
    @(44395, 0, 0, 0),
    @(44396, 0, 0, 0),
    @(44397, 0, 0, 0),
    @(44398, 0, 0, 0),
    @(44399, 0, 0, 0),
    @(44400, 0, 0, 0),
    @(44401, 0, 0, 0),
    @(44402, 1, 1, 26.76659528907923)
)

$startRow = 302
foreach ($row in $newData) {
    $r = $startRow
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $startRow++
}

# Mirror the date-column formatting (style index 2: bordered, centered, yyyy-mm-dd hh:mm:ss)
# from the last pre-existing row (A301) onto the newly appended date cells (A302:A328).
$ws.Range("A301").Copy() | Out-Null
$ws.Range("A302:A328").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
